# Added share price on 25/09/2017
# (commit message as supplied; actual content change is the new
# "Generics" rows appended to Sheet1's knowledge table)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Row 13: extend the "Generic Type" description and add an example ---
$ws1.Range("B13").Value = "A class or interface that is parameterized over types.`nProvides Stronger type checks at compile time.`nElimination of casts`nEnabling programmers to implement generic algorithms"
$ws1.Range("B13").Style = "Normal"
$ws1.Range("B13").VerticalAlignment = -4160
$ws1.Range("B13").WrapText = $true

$ws1.Range("C13").Value = "public class Box<T> {`n    // T stands for " + [char]34 + "Type" + [char]34 + "`n    private T t;`n    public void set(T t) { this.t = t; }`n    public T get() { return t; }`n}"
$ws1.Range("C13").VerticalAlignment = -4160
$ws1.Range("C13").WrapText = $true
$ws1.Rows.Item(13).RowHeight = 90

# --- Row 14: Type parameters ---
# (cell write order below matches the shared-string append order of the
#  target workbook: A14, B14, A15, C14, C15, A16, A17, C17, C16)
$ws1.Range("A14").Value = "Type parameters"
$ws1.Range("B14").Value = "E - Element (used extensively by the Java Collections Framework)`nK - Key`nN - Number`nT - Type`nV - Value`nS,U,V etc. - 2nd, 3rd, 4th types"
$ws1.Range("B14").VerticalAlignment = -4160
$ws1.Range("B14").WrapText = $true
$ws1.Range("E14").Value = "Generic Type"
$ws1.Rows.Item(14).RowHeight = 90

# --- Row 15: Type argument ---
$ws1.Range("A15").Value = "Type argument"
$ws1.Range("C14").Value = "Foo<T>"
$ws1.Range("C15").Value = "Foo<String>"

# --- Row 16: Diamond ---
$ws1.Range("A16").Value = "Diamond"

# --- Row 17: Multiple Type Parameters ---
$ws1.Range("A17").Value = "Multiple Type Parameters"
$ws1.Range("C17").Value = "public interface Pair<K, V> {`n    public K getKey();`n    public V getValue();`n}`n"
$ws1.Range("C17").VerticalAlignment = -4160
$ws1.Range("C17").WrapText = $true
$ws1.Rows.Item(17).RowHeight = 75

$ws1.Range("C16").Value = 'Box<Integer> integerBox = new Box<>();' + "`n" + 'OrderedPair<String, Integer> p1 = new OrderedPair<>("Even", 8);'
$ws1.Range("C16").VerticalAlignment = -4160
$ws1.Range("C16").WrapText = $true
$ws1.Rows.Item(16).RowHeight = 30

# --- Column widths: widen column C (code samples) and E (new label column) ---
$ws1.Columns.Item(3).ColumnWidth = 62
$ws1.Columns.Item(5).ColumnWidth = 11.67

# --- View state: Sheet1 becomes the active/selected tab (was Sheet2) ---
$ws1.Activate()
$ws1.Range("A16:XFD16").Select()
